# Generate Report for Handback
# Two handback md files got new UUIDs (one of them merged with a previously
# separate UUID), so every cell that referenced the old file names / xlf
# names / timestamps needs to be refreshed across the three sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Old identifiers
# ---------------------------------------------------------------------
$old1 = "6fcdd4e4-b5d7-4f1a-bbba-e5fa5aad9fc6"
$old2 = "97f69da1-c9e3-46fd-943a-dd10195ca5ab"

# New identifiers (old1 -> new1, old2 -> new2)
$new1 = "4d98daf5-3549-4f81-8d4b-c82bea531f69"
$new2 = "ffff3b4e6c6b-98f2-4091-9808-7fa0bc9e39d8"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value = ($new1 + ".md")
$ws.Range("B2").Value = ("e2e\" + $new1 + ".md")
$ws.Range("G2").Value = "2016-08-23 19:06:42"

$ws.Range("A3").Value = ($new2 + ".md")
$ws.Range("B3").Value = ("e2e\" + $new2 + ".md")
$ws.Range("G3").Value = "2016-08-23 19:06:42"

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A2").Value = ($new1 + ".md")
$ws.Range("G2").Value = ($new1 + ".77613b4ad2a8b46aea42ac927a4e36deca9a5c06.zh-cn.xlf")
$ws.Range("H2").Value = "2016-08-23 19:06:37"
$ws.Range("I2").Value = ($new1 + ".md")
$ws.Range("J2").Value = ($new1 + ".77613b4ad2a8b46aea42ac927a4e36deca9a5c06.zh-cn.xlf")
$ws.Range("K2").Value = "2016-08-23 19:06:53"

$ws.Range("A3").Value = ($new2 + ".md")
$ws.Range("G3").Value = ($new1 + ".77613b4ad2a8b46aea42ac927a4e36deca9a5c06.zh-cn.xlf")
$ws.Range("H3").Value = "2016-08-23 19:06:37"
$ws.Range("I3").Value = ($new2 + ".md")
$ws.Range("J3").Value = ($new1 + ".77613b4ad2a8b46aea42ac927a4e36deca9a5c06.zh-cn.xlf")
$ws.Range("K3").Value = "2016-08-23 19:06:53"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A2").Value = ($new1 + ".md")
$ws.Range("G2").Value = ($new1 + ".77613b4ad2a8b46aea42ac927a4e36deca9a5c06.de-de.xlf")
$ws.Range("H2").Value = "2016-08-23 19:06:42"
$ws.Range("I2").Value = ($new1 + ".md")
$ws.Range("J2").Value = ($new1 + ".77613b4ad2a8b46aea42ac927a4e36deca9a5c06.de-de.xlf")
$ws.Range("K2").Value = "2016-08-23 19:07:02"

$ws.Range("A3").Value = ($new2 + ".md")
$ws.Range("G3").Value = ($new1 + ".77613b4ad2a8b46aea42ac927a4e36deca9a5c06.de-de.xlf")
$ws.Range("H3").Value = "2016-08-23 19:06:42"
$ws.Range("I3").Value = ($new2 + ".md")
$ws.Range("J3").Value = ($new1 + ".77613b4ad2a8b46aea42ac927a4e36deca9a5c06.de-de.xlf")
$ws.Range("K3").Value = "2016-08-23 19:07:02"
